$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 3: 28/07/2022 -> 28-07-2022, D3 0->2, G3 0->2
Set-TextDate "A3" "28-07-2022"
$ws.Range("D3").Value = 2
$ws.Range("G3").Value = 2

# Row 4: 01/08/2022 -> 01-08-2022, D4 0->1, E4 0->1, H4 1->0
Set-TextDate "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: 04/08/2022 -> 04-08-2022, D5 0->1, E5 0->1, H5 1->0
Set-TextDate "A5" "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6: 08/08/2022 -> 08-08-2022 (no value changes)
Set-TextDate "A6" "08-08-2022"

# Row 7: 11/08/2022 -> 11-08-2022 (no value changes)
Set-TextDate "A7" "11-08-2022"

# Row 8: 15/08/2022 -> 15-08-2022 (no value changes)
Set-TextDate "A8" "15-08-2022"

# Row 9: 18/08/2022 -> 18-08-2022 (no value changes)
Set-TextDate "A9" "18-08-2022"

# Row 10: 22/08/2022 -> 22-08-2022 (no value changes)
Set-TextDate "A10" "22-08-2022"

# Row 11: 25/08/2022 -> 25-08-2022 (no value changes)
Set-TextDate "A11" "25-08-2022"

# Row 12: 29/08/2022 -> 29-08-2022 (no value changes)
Set-TextDate "A12" "29-08-2022"

# Row 13: 01/09/2022 -> 01-09-2022 (no value changes)
Set-TextDate "A13" "01-09-2022"

# Row 14: 05/09/2022 -> 05-09-2022 (no value changes)
Set-TextDate "A14" "05-09-2022"

# Row 15: 08/09/2022 -> 08-09-2022 (no value changes)
Set-TextDate "A15" "08-09-2022"

# Row 16: 12/09/2022 -> 12-09-2022 (no value changes)
Set-TextDate "A16" "12-09-2022"

# Row 17: 15/09/2022 -> 15-09-2022 (no value changes)
Set-TextDate "A17" "15-09-2022"

# Row 18: 19/09/2022 -> 19-09-2022 (no value changes)
Set-TextDate "A18" "19-09-2022"

# Row 19: 22/09/2022 -> 22-09-2022 (no value changes)
Set-TextDate "A19" "22-09-2022"

# Row 20: 26/09/2022 -> 26-09-2022 (no value changes)
Set-TextDate "A20" "26-09-2022"

# Row 21: 29/09/2022 -> 29-09-2022 (no value changes)
Set-TextDate "A21" "29-09-2022"
